$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row: new columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy style from an existing header cell (e.g. AB1) to the new headers
$ws.Range("AB1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wins = 67
$losses = 95
$ties = 0

for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 29).Value = $wins    # AC
    $ws.Cells.Item($r, 30).Value = $losses  # AD
    $ws.Cells.Item($r, 31).Value = $ties    # AE
}
